$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated cryptocurrency market data. Values that could be
# misread as numbers by Excel (e.g. "605.38", "1.00") are entered with a
# leading apostrophe so they stay plain text, matching the source data.

$ws.Range("D2").Value = "64.809.17"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "3.166.21"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'605.38"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "'145.58"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.158.18"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "'0.151"
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").Value = "'0.474"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "'35.73"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "3.681.53"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "64.679.82"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "3.182.07"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "'6.93"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "'484.31"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").Value = "'14.76"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").Value = "'0.718"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "'7.74"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("D24").Value = "'87.86"
$ws.Range("E24").Value = "  +4.90%  "
$ws.Range("D25").Value = "'13.60"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'2.79"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").Value = "'8.52"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "'7.21"
$ws.Range("E29").Value = "  +4.22%  "
$ws.Range("D30").Value = "'2.10"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "'0.113"
$ws.Range("E31").Value = "  -7.14%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'27.23"
$ws.Range("E32").Value = "  +3.64%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").Value = "'2.72"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").Value = "'1.11"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("D36").Value = "'6.10"
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("D37").Value = "0.0₃0771"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("D38").Value = "'53.05"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'3.07"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").Value = "'446.67"
$ws.Range("E40").Value = "  -2.57%  "
$ws.Range("D41").Value = "'0.0398"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").Value = "'0.120"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").Value = "'8.32"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").Value = "2.879.42"
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("D45").Value = "'0.265"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.50"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.27"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'26.33"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "'0.998"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "'121.87"
$ws.Range("E51").Value = "  +2.36%  "
